# "added update budget process"
#
# budget sheet: a new monthly "Food" budget entry is inserted at the top of
# the 45658 (period 1) group, shifting every later row down by one; the
# final (45748) period's "Transportation" budget amount changes from 250 to
# 300, and a brand-new "Utilities" row is appended for that same period.
#
# category sheet: the stray "Meow" category (row 10) is removed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# budget sheet
# ---------------------------------------------------------------------
$budget = $wb.Worksheets.Item("budget")
$dateFmt = "YYYY-MM-DD HH:MM:SS"

# Shift existing data rows 2..19 down to 3..20 (process bottom-up so we
# never overwrite a row before it has been read).
for ($r = 19; $r -ge 2; $r--) {
    $dst = $r + 1
    $budget.Cells.Item($dst, 1).Value = $budget.Cells.Item($r, 1).Value2
    $budget.Cells.Item($dst, 1).NumberFormat = $dateFmt
    $budget.Cells.Item($dst, 2).Value = $budget.Cells.Item($r, 2).Value2
    $budget.Cells.Item($dst, 3).Value = $budget.Cells.Item($r, 3).Value2
    $budget.Cells.Item($dst, 4).Value = $budget.Cells.Item($r, 4).Value2
    $budget.Cells.Item($dst, 5).Value = $budget.Cells.Item($r, 5).Value2
}

# New row 2: Food budget for the first (45658) period.
$budget.Cells.Item(2, 1).Value = 45658
$budget.Cells.Item(2, 1).NumberFormat = $dateFmt
$budget.Cells.Item(2, 2).Value = "Food"
$budget.Cells.Item(2, 3).Value = 600
$budget.Cells.Item(2, 4).Value = 1
$budget.Cells.Item(2, 5).Value = 2

# Row 20 (old row 19 shifted down) is the 45748/Transportation row whose
# monthly_budget increases from 250 to 300.
$budget.Cells.Item(20, 3).Value = 300

# New row 21: Utilities budget for the last (45748) period.
$budget.Cells.Item(21, 1).Value = 45748
$budget.Cells.Item(21, 1).NumberFormat = $dateFmt
$budget.Cells.Item(21, 2).Value = "Utilities"
$budget.Cells.Item(21, 3).Value = 50
$budget.Cells.Item(21, 4).Value = 20
$budget.Cells.Item(21, 5).Value = 5

# ---------------------------------------------------------------------
# category sheet
# ---------------------------------------------------------------------
$category = $wb.Worksheets.Item("category")
$category.Rows.Item(10).Delete()
